# kwitansi_dinas.docx update: add the "SISA ANGGARAN" row label.
#
# The second table in the document (the SPJ/receipt detail grid whose
# first cell reads "SKPD") has a totals section near the bottom:
#   row "JUMLAH"              -> subtotal row
#   row (currently blank)     -> should read "SISA ANGGARAN"
#   row "${sisa_sesudah}"     -> remaining-balance values
#
# Row 15's first (merged, gridSpan=5) cell is that blank label cell; it
# needs a bold, 9pt ("half-points" 18) run reading "SISA ANGGARAN",
# mirroring the "JUMLAH" row directly above it.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)

$labelCell = $tbl.Cell(15, 1)
$labelCell.Range.Paragraphs.Item(1).Range.Text = "SISA ANGGARAN"

# Re-fetch the range covering just the new run (not the whole cell, so
# the paragraph mark itself is left untouched) and apply bold + 9pt,
# forcing both the Latin and complex-script size so the run explicitly
# carries <w:sz .../><w:szCs .../> like its sibling rows do.
$runRange = $tbl.Cell(15, 1).Range.Paragraphs.Item(1).Range
$runRange.Font.Bold = $true
$runRange.Font.Size = 9
$runRange.Font.SizeBi = 9

# Note: the other hunk in the source diff only bumps the internal
# o:OLEObject/@ObjectID bookkeeping id on the header logo (same
# ProgID/ShapeID/r:id - no content/visual change). That attribute is an
# implementation-detail counter Word mints/rewrites internally when it
# resaves a document containing an embedded OLE object; it is not
# exposed anywhere in the Word COM object model (no OLEFormat/Field/
# Shape property maps to it), so it cannot be targeted from automation
# code and is intentionally left as-is.
